$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the leading "Unnamed"/empty columns (old A:D), shifting
# "Mem ID"/"NAME"/"LAST NAME"/"ID"/"NUMBER"/"DATE" (old E:J) left into A:F.
$ws.Range("A1:D3").Delete(-4159)

# The data rows must keep their literal text (e.g. leading zeros, no date
# parsing), so format them as text before writing the new values, then
# drop back to the default (unstyled) look the diff expects.
$data = $ws.Range("A2:F3")
$data.ClearFormats()
$data.NumberFormat = "@"

$ws.Range("A2").Value = "507365"
$ws.Range("B2").Value = "ilia"
$ws.Range("C2").Value = "valaei"
$ws.Range("D2").Value = "0441201425"
$ws.Range("E2").Value = "09307637377"
$ws.Range("F2").Value = "02/20/2023"

$ws.Range("A3").Value = "557690"
$ws.Range("B3").Value = "sdfgsgsdf"
$ws.Range("C3").Value = "gsdgsdgsdgdsg"
$ws.Range("D3").Value = "gsdgsdgd"
$ws.Range("E3").Value = "ssdgsdgdsg"
$ws.Range("F3").Value = "02/20/2023"

$data.ClearFormats()
